# "Added filter by geo on collecting cript"
# Adds a new column H "Used for RQ" that tags each indicator row with the
# research question(s) it supports, and fixes a stray formatting outlier
# on A11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column H: width + header + data, formats copied from existing
#     cells so we reuse the workbook's existing style entries instead of
#     minting new ones.
$ws.Columns.Item(8).ColumnWidth = 15.33

$ws.Range("A4").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A5").Copy()
$ws.Range("H5:H19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Values (filled in the order that reproduces the author's shared
#     string table: header, then the "core" GDP row, then the rest).
$ws.Range("H4").Value = "Used for RQ"
$ws.Range("H11").Value = "RQ1"
$ws.Range("H5").Value = "RQ3"
$ws.Range("H7").Value = "RQ1, RQ3"
$ws.Range("H6").Value = "RQ2"
$ws.Range("H8").Value = "RQ1, RQ3"
$ws.Range("H9").Value = "RQ1, RQ3"
$ws.Range("H10").Value = "RQ1, RQ3"
$ws.Range("H12").Value = "RQ2"
$ws.Range("H13").Value = "RQ1, RQ3"
$ws.Range("H14").Value = "RQ1"
$ws.Range("H15").Value = "RQ2"
$ws.Range("H16").Value = "RQ1, RQ3"
$ws.Range("H17").Value = "RQ3"
$ws.Range("H18").Value = "RQ2"
$ws.Range("H19").Value = "RQ2"

# --- Row heights: header row + the rows whose H-text needs more wrap
#     room grew taller.
$ws.Rows.Item(4).RowHeight = 31.5
$ws.Rows.Item(8).RowHeight = 47.25
$ws.Rows.Item(9).RowHeight = 47.25
$ws.Rows.Item(10).RowHeight = 47.25
$ws.Rows.Item(11).RowHeight = 31.5

# --- A11 had a stray one-off style (border only, no font) left over from
#     when that row was pasted in; realign it with the rest of column A.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A11").Value = 7

# --- Selection / scroll state to match the saved view.
$ws.Range("I13").Select()
